$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 556; existing rows 556-628 shift down to 557-629.
$ws.Rows(556).Insert()

# Populate the newly inserted row 556 with its values.
$ws.Range("A556").Value = 4
$ws.Range("B556").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C556").Value = "Los Lagos"
$ws.Range("D556").Value = 44776
$ws.Range("E556").Value = 10
$ws.Range("F556").Value = 100112004
$ws.Range("G556").Value = "Cebolla"
$ws.Range("H556").Value = "Sin especificar"
$ws.Range("I556").Value = "2a (guarda)"
$ws.Range("J556").Value = 150
$ws.Range("K556").Value = 8500
$ws.Range("L556").Value = 8500
$ws.Range("M556").Value = 8500
$ws.Range("N556").Value = "$/malla 18 kilos"
$ws.Range("O556").Value = "Región de O'Higgins"
$ws.Range("P556").Value = 472
$ws.Range("Q556").Value = 18
$ws.Range("R556").Value = "Hortaliza"
